$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.83
$ws.Range("I2").Value = 3.8
$ws.Range("L2").Value = 4.33
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 12
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2
$ws.Range("U2").Value = 1.67
$ws.Range("V2").Value = 2.1
$ws.Range("X2").Value = 9.5
$ws.Range("Z2").Value = 17
$ws.Range("AC2").Value = 12
$ws.Range("AJ2").Value = 13
$ws.Range("AL2").Value = 29
$ws.Range("AM2").Value = 34
$ws.Range("AO2").Value = 10
$ws.Range("AY2").Value = 26
$ws.Range("BB2").Value = 151
$ws.Range("K3").Value = 2.25
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 3.5
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.93
$ws.Range("AL3").Value = 34
$ws.Range("AN3").Value = 3.75
$ws.Range("G4").Value = 2.55
$ws.Range("I4").Value = 2.45
$ws.Range("J4").Value = 3.2
$ws.Range("K4").Value = 2.3
$ws.Range("O4").Value = 1.22
$ws.Range("P4").Value = 3.85
$ws.Range("Q4").Value = 1.67
$ws.Range("R4").Value = 2.15
$ws.Range("U4").Value = 1.55
$ws.Range("V4").Value = 2.3
$ws.Range("AG4").Value = 126
$ws.Range("J5").Value = 4.5
$ws.Range("K5").Value = 2.1
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 3.4
$ws.Range("Q5").Value = 2.03
$ws.Range("R5").Value = 1.83
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63
$ws.Range("W5").Value = 11
$ws.Range("Y5").Value = 15
$ws.Range("AC5").Value = 9.5
$ws.Range("AH5").Value = 7
$ws.Range("AM5").Value = 29
$ws.Range("AO5").Value = 23
$ws.Range("AQ5").Value = 81
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 2.63
$ws.Range("AW5").Value = 3.75
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 1.9
$ws.Range("J7").Value = 4.75
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 2.63
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 2.75
$ws.Range("Q7").Value = 2.25
$ws.Range("R7").Value = 1.62
$ws.Range("S7").Value = 1.5
$ws.Range("T7").Value = 2.5
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.73
$ws.Range("W7").Value = 10
$ws.Range("X7").Value = 21
$ws.Range("Z7").Value = 41
$ws.Range("AA7").Value = 41
$ws.Range("AB7").Value = 41
$ws.Range("AC7").Value = 8
$ws.Range("AE7").Value = 17
$ws.Range("AF7").Value = 67
$ws.Range("AG7").Value = 1000
$ws.Range("AH7").Value = 6
$ws.Range("AI7").Value = 8.5
$ws.Range("AJ7").Value = 9
$ws.Range("AK7").Value = 15
$ws.Range("AL7").Value = 17
$ws.Range("AM7").Value = 34
$ws.Range("AN7").Value = 6
$ws.Range("AO7").Value = 23
$ws.Range("AP7").Value = 34
$ws.Range("AQ7").Value = 81
$ws.Range("AR7").Value = 126
$ws.Range("AS7").Value = 301
$ws.Range("AT7").Value = 2.5
$ws.Range("AU7").Value = 9
$ws.Range("AV7").Value = 67
$ws.Range("AW7").Value = 3.75
$ws.Range("AX7").Value = 11
$ws.Range("AY7").Value = 23
$ws.Range("AZ7").Value = 41
$ws.Range("BA7").Value = 67
$ws.Range("BB7").Value = 201
$ws.Range("G8").Value = 3.8
$ws.Range("H8").Value = 3.6
$ws.Range("I8").Value = 1.85
$ws.Range("J8").Value = 4.5
$ws.Range("K8").Value = 2.2
$ws.Range("L8").Value = 2.5
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("O8").Value = 1.29
$ws.Range("P8").Value = 3.5
$ws.Range("Q8").Value = 1.95
$ws.Range("R8").Value = 1.85
$ws.Range("S8").Value = 1.4
$ws.Range("T8").Value = 2.75
$ws.Range("U8").Value = 1.8
$ws.Range("V8").Value = 1.91
$ws.Range("W8").Value = 11
$ws.Range("X8").Value = 19
$ws.Range("Y8").Value = 13
$ws.Range("Z8").Value = 41
$ws.Range("AA8").Value = 34
$ws.Range("AB8").Value = 41
$ws.Range("AC8").Value = 11
$ws.Range("AD8").Value = 7
$ws.Range("AE8").Value = 15
$ws.Range("AF8").Value = 51
$ws.Range("AG8").Value = 251
$ws.Range("AH8").Value = 7.5
$ws.Range("AI8").Value = 9
$ws.Range("AJ8").Value = 8.5
$ws.Range("AK8").Value = 15
$ws.Range("AL8").Value = 15
$ws.Range("AM8").Value = 26
$ws.Range("AN8").Value = 6
$ws.Range("AO8").Value = 21
$ws.Range("AP8").Value = 29
$ws.Range("AQ8").Value = 67
$ws.Range("AR8").Value = 101
$ws.Range("AS8").Value = 201
$ws.Range("AT8").Value = 2.75
$ws.Range("AU8").Value = 8
$ws.Range("AV8").Value = 51
$ws.Range("AW8").Value = 4
$ws.Range("AX8").Value = 10
$ws.Range("AY8").Value = 21
$ws.Range("AZ8").Value = 34
$ws.Range("BA8").Value = 51
$ws.Range("BB8").Value = 151
$ws.Range("G9").Value = 1.75
$ws.Range("H9").Value = 3.4
$ws.Range("J9").Value = 2.5
$ws.Range("K9").Value = 2.05
$ws.Range("L9").Value = 5
$ws.Range("N9").Value = 8
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 3
$ws.Range("R9").Value = 1.67
$ws.Range("T9").Value = 2.63
$ws.Range("V9").Value = 1.73
$ws.Range("W9").Value = 6
$ws.Range("X9").Value = 7.5
$ws.Range("Z9").Value = 13
$ws.Range("AB9").Value = 34
$ws.Range("AC9").Value = 8
$ws.Range("AD9").Value = 6.5
$ws.Range("AE9").Value = 19
$ws.Range("AF9").Value = 67
$ws.Range("AH9").Value = 11
$ws.Range("AI9").Value = 23
$ws.Range("AJ9").Value = 15
$ws.Range("AK9").Value = 51
$ws.Range("AL9").Value = 41
$ws.Range("AM9").Value = 41
$ws.Range("AN9").Value = 3.6
$ws.Range("AO9").Value = 9.5
$ws.Range("AP9").Value = 23
$ws.Range("AQ9").Value = 34
$ws.Range("AR9").Value = 51
$ws.Range("AS9").Value = 201
$ws.Range("AT9").Value = 2.63
$ws.Range("AU9").Value = 9
$ws.Range("AV9").Value = 67
$ws.Range("AW9").Value = 6.5
$ws.Range("AY9").Value = 34
$ws.Range("AZ9").Value = 101
$ws.Range("BA9").Value = 126
$ws.Range("BB9").Value = 301
$ws.Range("G10").Value = 1.85
$ws.Range("I10").Value = 3.7
$ws.Range("M10").Value = 1.05
$ws.Range("N10").Value = 8.5
$ws.Range("O10").Value = 1.29
$ws.Range("P10").Value = 3.5
$ws.Range("Q10").Value = 1.95
$ws.Range("R10").Value = 1.85
$ws.Range("S10").Value = 1.4
$ws.Range("T10").Value = 2.75
$ws.Range("X10").Value = 9
$ws.Range("AC10").Value = 9.5
$ws.Range("AI10").Value = 19
$ws.Range("AJ10").Value = 13
$ws.Range("AN10").Value = 4
$ws.Range("AT10").Value = 2.75
$ws.Range("AX10").Value = 21
$ws.Range("AY10").Value = 29
$ws.Range("AZ10").Value = 67
$ws.Range("BB10").Value = 400
$ws.Range("H12").Value = 3.2
$ws.Range("K12").Value = 2.1
$ws.Range("L12").Value = 3.4
$ws.Range("Q12").Value = 1.9
$ws.Range("R12").Value = 1.9
$ws.Range("S12").Value = 1.4
$ws.Range("T12").Value = 2.75
$ws.Range("U12").Value = 1.73
$ws.Range("V12").Value = 2
$ws.Range("W12").Value = 8.5
$ws.Range("AB12").Value = 29
$ws.Range("AC12").Value = 10
$ws.Range("AL12").Value = 23
$ws.Range("AM12").Value = 34
$ws.Range("AP12").Value = 23
$ws.Range("AR12").Value = 67
$ws.Range("AT12").Value = 2.75
$ws.Range("AU12").Value = 8
$ws.Range("AY12").Value = 26
$ws.Range("G13").Value = 3
$ws.Range("I13").Value = 2.1
$ws.Range("L13").Value = 3
$ws.Range("O13").Value = 1.33
$ws.Range("P13").Value = 3.25
$ws.Range("Q13").Value = 2.08
$ws.Range("R13").Value = 1.73
$ws.Range("U13").Value = 1.91
$ws.Range("V13").Value = 1.91
$ws.Range("W13").Value = 8.5
$ws.Range("Y13").Value = 11
$ws.Range("AA13").Value = 26
$ws.Range("AB13").Value = 34
$ws.Range("AG13").Value = 301
$ws.Range("AJ13").Value = 9.5
$ws.Range("AK13").Value = 21
$ws.Range("AW13").Value = 4.33
$ws.Range("AX13").Value = 13
